# CM70_TestData_ManuallyCreateBankAccountStatement_21C.xlsx
# "Add files via upload" / "Anu - Cash Management Files Uploaded"
#
# The previous upload had the Input_Value sheet pre-populated with a live
# Oracle Cloud URL + implementation-user credentials in Y2:AA2 (with Y2
# hyperlinked to that URL). This revision blanks those three cells back
# out (and drops the now-unused hyperlink) so the sheet is shipped as a
# clean, credential-free template; the analyst's selection is left
# sitting on that same Y2:AA2 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

$ws.Activate()

# Drop the hyperlink that lived on Y2 (https://edrx.fa.us2.oraclecloud.com)
# before clearing the cell text, then blank out the URL/user/password
# triple in Y2:AA2 - content only, formatting (styles) stays untouched.
$ws.Range("Y2").Hyperlinks.Delete()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()

# Leave the selection on the now-empty credentials block.
$ws.Range("Y2:AA2").Select()
